$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select rows 3-7 (the CLICK_JS/WAIT block for the EnableCertificate steps),
# then delete them - this shifts the remaining rows (old 8-11, i.e. the
# CLICK/ENTERTEXT/WAIT/VERIFY_WEBELEMENT_PRESENT block) up to become rows 3-6.
$ws.Rows("3:7").Select()
$excel.Selection.Delete()
